# Apply "Notes and examples variables" update to the training-schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 -------------------------------------------------------------
$ws.Range("D10").Value = "day-6"
$ws.Range("E10").Value = "3/30/2023"
$ws.Range("F10").Value = "js advanced "
$ws.Range("G10").Value = "09:00PM-10:00PM"

# --- Row 11 -------------------------------------------------------------
$ws.Range("E11").Value = "3/31/2023"
$ws.Range("F11").Value = "js validations,"
$ws.Range("G11").Value = "09:00PM-10:00PM"

# --- Row 12 -------------------------------------------------------------
$ws.Range("E12").Value = "4/3/2023"
$ws.Range("F12").Value = "Bootstrap"
$ws.Range("G12").Value = "9:30PM-10:20PM"

# --- Row 13 -------------------------------------------------------------
$ws.Range("E13").Value = "4/6/2023"
$ws.Range("F13").Value = "Typescript"
$ws.Range("G13").Value = "9:00pm-10:00pm"

# --- Row 14 -------------------------------------------------------------
$ws.Range("E14").Value = "4/7/2023"
$ws.Range("F14").Value = "Angular intro"
$ws.Range("G14").Value = "9:00pm-10:00pm"

# --- Row 15 -------------------------------------------------------------
$ws.Range("E15").Value = "4/10/2023"
$ws.Range("F15").Value = "Angular components,commands ,basic structure"
$ws.Range("G15").Value = "9:05-10:00pm"

# --- Row 16 -------------------------------------------------------------
$ws.Range("E16").Value = "4/11/2023"
$ws.Range("F16").Value = "routing ,directives,data binding …."
$ws.Range("G16").Value = "9:00-10:05pm"
$ws.Range("H16").Value = "direcgtives revise "

# --- Row 17 -------------------------------------------------------------
$ws.Range("E17").Value = "4/12/2023"
$ws.Range("F17").Value = "pipes,forms "

# --- Row 18 -------------------------------------------------------------
$ws.Range("E18").Value = "4/19/2023"
$ws.Range("F18").Value = "authguards"

# --- Row 19 -------------------------------------------------------------
$ws.Range("E19").Value = "4/24/2023"
$ws.Range("F19").Value = "Json files,services,crud operations json…,httpclient"
$ws.Range("G19").Value = "9:45PM-10:25Pm"

# --- Row 20 (Timings/Tasks entered before Concepts, matching author order) --
$ws.Range("E20").Value = "4/26/2023"
$ws.Range("G20").Value = "10:30PM -11:30PM"
$ws.Range("H20").Value = "complete java setup  and eclipse "
$ws.Range("F20").Value = "Java Intro,java setup"

# --- Row 21 -------------------------------------------------------------
$ws.Range("E21").Value = "4/28/2023"
$ws.Range("F21").Value = "execution flow,basic variables,identifiers"
$ws.Range("G21").Value = "10:30PM -11:25PM"

# --- Row 22 -------------------------------------------------------------
$ws.Range("E22").Value = "5/2/2023"
$ws.Range("F22").Value = "Variables,datatypes,Eclipse"
$ws.Range("G22").Value = "10:30PM -11:25PM"

# --- Row 23 -------------------------------------------------------------
$ws.Range("E23").Value = "5/5/2023"
$ws.Range("F23").Value = "Variables instance,static"
$ws.Range("G23").Value = "10:30PM -11:30PM"

# --- Fix number-formats / borders that drifted from plain value entry ----
# G10, G11 need the header-row time style (same as G5:G9)
$ws.Range("G9").Copy() | Out-Null
$ws.Range("G10:G11").PasteSpecial(-4122) | Out-Null

# E14 should use the short-date style used by E8:E13 (not the E15:E21 style)
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

# G14 and G21 drop back to the plain style used elsewhere in column G
$ws.Range("G24").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G21").PasteSpecial(-4122) | Out-Null

# E22 and E23 pick up the E15:E21 date style
$ws.Range("E21").Copy() | Out-Null
$ws.Range("E22:E23").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column H is widened to fit the new "Tasks" text ---------------------
# (the engine rounds ColumnWidth to 1/6-character steps when it persists the
#  XML <col> width, so the input is pre-compensated to land on 32.33203125)
$ws.Columns.Item(8).ColumnWidth = 31.4986979166667

# --- Restore the view/selection state recorded in the saved file ---------
$ws.Range("H23").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
